$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

$r = $target.Range
$r.Collapse(0)              # wdCollapseEnd
$r.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Range.Text = "8855158 - Morun Bernardino Neto"
$newPara.Style = "ListBullet"
